$d = $word.ActiveDocument

$replacements = @(
    @("2025-06-16 Monday", "2025-06-17 Tuesday"),
    @("13-9=", "76-40="),
    @("47-22=", "51-40="),
    @("52-50=", "45+4="),
    @("85-66=", "5+9="),
    @("77-12=", "62-23="),
    @("45+53=", "66-38="),
    @("74-52=", "63+2="),
    @("15+79=", "20-0="),
    @("48-46=", "91-17="),
    @("72-65=", "59-24="),
    @("56-3=", "74-26="),
    @("85-4=", "23+38="),
    @("17+63=", "62+20="),
    @("63+3=", "13+69="),
    @("37-23=", "24+20="),
    @("80-62=", "61-7="),
    @("17+78=", "90-18="),
    @("16+57=", "57-31="),
    @("0+20=", "50-38="),
    @("25+31=", "63+35="),
    @("29+10=", "39+48="),
    @("17-16=", "67-1="),
    @("54-21=", "2+35="),
    @("53-49=", "38+3="),
    @("33-2=", "10+78="),
    @("51+36=", "37+11="),
    @("62-21=", "89-3="),
    @("10+67=", "95-83="),
    @("6+84=", "80-12="),
    @("52-16=", "49-1="),
    @("17+49=", "39-29="),
    @("3+61=", "97-76="),
    @("2+8=", "33+48="),
    @("51+10=", "30-5="),
    @("90-62=", "40+16="),
    @("0+38=", "88-30="),
    @("30+6=", "41-13="),
    @("75-31=", "46-36="),
    @("24+36=", "71-3="),
    @("21+37=", "0+27="),
    @("41+9=", "54+38="),
    @("52-4=", "71-15="),
    @("49+42=", "92+3="),
    @("27+64=", "20+73="),
    @("37+6=", "78+10="),
    @("58+1=", "58+38="),
    @("28+32=", "10+70="),
    @("95-87=", "29+53="),
    @("23+23=", "38-11="),
    @("68-53=", "54+26="),
    @("86-81=", "26+38="),
    @("51+13=", "45-7="),
    @("61+10=", "50-34="),
    @("96-14=", "94-65="),
    @("84-69=", "28+62="),
    @("61+7=", "44-27="),
    @("80-34=", "8+40="),
    @("68-10=", "40+14="),
    @("11-5=", "97-65="),
    @("78-51=", "99-85="),
    @("73-58=", "3+35="),
    @("20-11=", "46-33="),
    @("74-39=", "84-20="),
    @("47+16=", "50-8="),
    @("94-45=", "75+9="),
    @("85-42=", "70+25="),
    @("33+5=", "21+36="),
    @("63-32=", "36-29="),
    @("23-23=", "27-7="),
    @("13+71=", "34-9="),
    @("63+17=", "83-26="),
    @("63-20=", "50+31="),
    @("84-79=", "52+9="),
    @("13+14=", "10+64="),
    @("66-17=", "27-3="),
    @("43+39=", "87-20="),
    @("86+5=", "66-59="),
    @("4+41=", "17+12="),
    @("37+5=", "57-0="),
    @("11-2=", "36-8="),
    @("66-55=", "98-47="),
    @("22+3=", "91-27="),
    @("1+58=", "89+1="),
    @("88-24=", "13+59="),
    @("87-85=", "62+15="),
    @("19+11=", "58-19="),
    @("14-10=", "70-8="),
    @("5+13=", "32+29="),
    @("60-2=", "49+24="),
    @("45+36=", "85-17="),
    @("7+59=", "4+29="),
    @("71-31=", "34-7="),
    @("26+4=", "89-62="),
    @("40+1=", "56-35="),
    @("25+12=", "44-38="),
    @("38+10=", "54-2="),
    @("45+31=", "32+16="),
    @("58+6=", "94-52="),
    @("67-63=", "96-67="),
    @("76-29=", "64+23="),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}

Write-Output "Done"
